# Rename the "wt" and "dcin5" sheets to include the "_log2_expression" suffix,
# matching the updated input-file naming convention, and move the active tab
# selection from "optimization_parameters" over to the newly renamed
# "dcin5_log2_expression" sheet.
$wb = $excel.ActiveWorkbook

$wsWt = $wb.Worksheets.Item("wt")
$wsWt.Name = "wt_log2_expression"

$wsDcin5 = $wb.Worksheets.Item("dcin5")
$wsDcin5.Name = "dcin5_log2_expression"

$wsDcin5.Activate()
